# Update "想去人数" (interested-count) figures on the "展览" and "全部类型" sheets
# F4: 7812 -> 7813
# F5: 5707 -> 5710

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 7813
    $ws.Range("F5").Value = 5710
}
